$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 64
$ws.Range("I64").Value = 3525.125
$ws.Range("J64").Value = 4700
$ws.Range("K64").Value = 3525.125
$ws.Range("L64").Value = 4700
$ws.Range("M64").Value = -3277.125
$ws.Range("N64").Value = -5196
# Row 67
$ws.Range("I67").Value = 3525.125
$ws.Range("J67").Value = 4700
$ws.Range("K67").Value = 3525.125
$ws.Range("L67").Value = 4700
$ws.Range("M67").Value = -2667.125
$ws.Range("N67").Value = -6416
# Row 76
$ws.Range("H76").Value = 3873.913
$ws.Range("I76").Value = 3687.9412
$ws.Range("J76").Value = 4400.8335
$ws.Range("K76").Value = 3687.9412
$ws.Range("L76").Value = 4400.8335
$ws.Range("M76").Value = -3372.9412
$ws.Range("N76").Value = -5030.8335
# Row 79
$ws.Range("H79").Value = 3873.913
$ws.Range("I79").Value = 3687.9412
$ws.Range("J79").Value = 4400.8335
$ws.Range("K79").Value = 3687.9412
$ws.Range("L79").Value = 4400.8335
$ws.Range("M79").Value = -2595.9412
$ws.Range("N79").Value = -6584.8335
# Row 132
$ws.Range("H132").Value = 1462.7963
$ws.Range("I132").Value = 1305.9387
$ws.Range("K132").Value = 3917.8161
$ws.Range("M132").Value = -1387.8161
# Row 137
$ws.Range("H137").Value = 3391.7856
$ws.Range("I137").Value = 1784.6086
$ws.Range("J137").Value = 5337.316
$ws.Range("K137").Value = 5353.825800000001
$ws.Range("L137").Value = 16011.948
$ws.Range("M137").Value = -2803.825800000001
$ws.Range("N137").Value = -21111.948

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1943.7333
$ws.Range("I2").Value = 2058.5833
$ws.Range("J2").Value = 1484.3334
$ws.Range("K2").Value = 2058.5833
$ws.Range("L2").Value = 1484.3334
$ws.Range("M2").Value = -1945.5833
$ws.Range("N2").Value = -1710.3334
# Row 63
$ws.Range("H63").Value = 3192.6155
$ws.Range("I63").Value = 1899.6
$ws.Range("K63").Value = 1899.6
$ws.Range("M63").Value = -1213.6
# Row 66
$ws.Range("H66").Value = 3192.6155
$ws.Range("I66").Value = 1899.6
$ws.Range("K66").Value = 9498
$ws.Range("M66").Value = -6066
# Row 116
$ws.Range("H116").Value = 1943.7333
$ws.Range("I116").Value = 2058.5833
$ws.Range("J116").Value = 1484.3334
$ws.Range("K116").Value = 2058.5833
$ws.Range("L116").Value = 1484.3334
$ws.Range("M116").Value = 235.4167000000002
$ws.Range("N116").Value = -6072.3334
# Row 123
$ws.Range("H123").Value = 51414.5
$ws.Range("J123").Value = 51414.5
$ws.Range("L123").Value = 51414.5
$ws.Range("N123").Value = -61214.5

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1943.7333
$ws.Range("I3").Value = 2058.5833
$ws.Range("J3").Value = 1484.3334
$ws.Range("K3").Value = 2058.5833
$ws.Range("L3").Value = 1484.3334
$ws.Range("M3").Value = -1944.5833
$ws.Range("N3").Value = -1712.3334
# Row 9
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()
# Row 26
$ws.Range("H26").Value = 16804.334
$ws.Range("I26").Value = 16804.334
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 16804.334
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -16512.334
$ws.Range("N26").ClearContents()
# Row 40
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
# Row 105
$ws.Range("H105").Value = 5281.4146
$ws.Range("I105").Value = 4923.6294
$ws.Range("J105").Value = 5971.4287
$ws.Range("K105").Value = 4923.6294
$ws.Range("L105").Value = 5971.4287
$ws.Range("M105").Value = -3176.6294
$ws.Range("N105").Value = -9465.4287
# Row 109
$ws.Range("H109").Value = 47684
$ws.Range("J109").Value = 47684
$ws.Range("L109").Value = 47684
$ws.Range("N109").Value = -50458

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2945.8918
$ws.Range("I31").Value = 2199.7856
$ws.Range("J31").Value = 5267.1113
$ws.Range("K31").Value = 2199.7856
$ws.Range("L31").Value = 5267.1113
$ws.Range("M31").Value = -1904.7856
$ws.Range("N31").Value = -5857.1113
# Row 34
$ws.Range("H34").Value = 2945.8918
$ws.Range("I34").Value = 2199.7856
$ws.Range("J34").Value = 5267.1113
$ws.Range("K34").Value = 2199.7856
$ws.Range("L34").Value = 5267.1113
$ws.Range("M34").Value = -1997.7856
$ws.Range("N34").Value = -5671.1113
# Row 122
$ws.Range("H122").Value = 11925.895
$ws.Range("I122").Value = 6309.8887
$ws.Range("J122").Value = 16980.3
$ws.Range("K122").Value = 18929.6661
$ws.Range("L122").Value = 50940.89999999999
$ws.Range("M122").Value = -16479.6661
$ws.Range("N122").Value = -55840.89999999999

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 5751461.5
$ws.Range("I5").Value = 657
$ws.Range("J5").Value = 16677990
$ws.Range("K5").Value = 1971
$ws.Range("L5").Value = 50033970
$ws.Range("M5").Value = -1859
$ws.Range("N5").Value = -50034194
# Row 12
$ws.Range("H12").Value = 26316032
$ws.Range("I12").Value = 76923310
$ws.Range("J12").Value = 243.92
$ws.Range("K12").Value = 230769930
$ws.Range("L12").Value = 731.76
$ws.Range("M12").Value = -230769757
$ws.Range("N12").Value = -1077.76
# Row 107
$ws.Range("H107").Value = 876.25
$ws.Range("I107").Value = 320.58334
$ws.Range("J107").Value = 1061.4722
$ws.Range("K107").Value = 961.7500200000001
$ws.Range("L107").Value = 3184.4166
$ws.Range("M107").Value = 958.2499799999999
$ws.Range("N107").Value = -7024.4166
# Row 135
$ws.Range("H135").Value = 5751461.5
$ws.Range("I135").Value = 657
$ws.Range("J135").Value = 16677990
$ws.Range("K135").Value = 5913
$ws.Range("L135").Value = 150101910
$ws.Range("M135").Value = -3378
$ws.Range("N135").Value = -150106980
# Row 136
$ws.Range("H136").Value = 2284.88
$ws.Range("I136").Value = 917.61536
$ws.Range("J136").Value = 3766.0833
$ws.Range("K136").Value = 2752.84608
$ws.Range("L136").Value = 11298.2499
$ws.Range("M136").Value = 2347.15392
$ws.Range("N136").Value = -21498.2499

$ws = $wb.Worksheets.Item("GSM")
# Row 19
$ws.Range("H19").Value = 5000
$ws.Range("J19").Value = 5000
$ws.Range("L19").Value = 5000
$ws.Range("N19").Value = -5576
# Row 70
$ws.Range("H70").Value = 4972.623
$ws.Range("I70").Value = 4727
$ws.Range("J70").Value = 5075.442
$ws.Range("K70").Value = 4727
$ws.Range("L70").Value = 5075.442
$ws.Range("M70").Value = -4457
$ws.Range("N70").Value = -5615.442
# Row 73
$ws.Range("H73").Value = 4972.623
$ws.Range("I73").Value = 4727
$ws.Range("J73").Value = 5075.442
$ws.Range("K73").Value = 4727
$ws.Range("L73").Value = 5075.442
$ws.Range("M73").Value = -3791
$ws.Range("N73").Value = -6947.442
# Row 123
$ws.Range("H123").Value = 28316.857
$ws.Range("J123").Value = 28316.857
$ws.Range("L123").Value = 28316.857
$ws.Range("N123").Value = -33216.857
# Row 132
$ws.Range("H132").Value = 8100.1763
$ws.Range("I132").Value = 21992.6
$ws.Range("K132").Value = 65977.79999999999
$ws.Range("M132").Value = -63447.79999999999

$ws = $wb.Worksheets.Item("LTW")
# Row 4
$ws.Range("H4").Value = 22300
$ws.Range("J4").Value = 22300
$ws.Range("L4").Value = 22300
$ws.Range("N4").Value = -22526
# Row 28
$ws.Range("H28").Value = 22300
$ws.Range("J28").Value = 22300
$ws.Range("L28").Value = 22300
$ws.Range("N28").Value = -22764
# Row 37
$ws.Range("H37").Value = 22300
$ws.Range("J37").Value = 22300
$ws.Range("L37").Value = 22300
$ws.Range("N37").Value = -22514
# Row 61
$ws.Range("H61").Value = 1696755.6
$ws.Range("I61").Value = 2438891
$ws.Range("J61").Value = 26951.25
$ws.Range("K61").Value = 2438891
$ws.Range("L61").Value = 26951.25
$ws.Range("M61").Value = -2438689
$ws.Range("N61").Value = -27355.25
# Row 113
$ws.Range("H113").Value = 1696755.6
$ws.Range("I113").Value = 2438891
$ws.Range("J113").Value = 26951.25
$ws.Range("K113").Value = 2438891
$ws.Range("L113").Value = 26951.25
$ws.Range("M113").Value = -2436721
$ws.Range("N113").Value = -31291.25
# Row 140
$ws.Range("H140").Value = 69834.73
$ws.Range("J140").Value = 69834.73
$ws.Range("L140").Value = 69834.73
$ws.Range("N140").Value = -80194.73

$ws = $wb.Worksheets.Item("WVR")
# Row 6
$ws.Range("H6").Value = 3381.2
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 3381.2
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 3381.2
$ws.Range("M6").ClearContents()
$ws.Range("N6").Value = -3611.2
